$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Version cell from "0.1" to "1.0" (force text so Excel doesn't coerce it to a number)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.0"

# Swap the step-2 content between TC2 (row 20) and TC3 (row 28):
# TC2's step-2 now describes cancelling a diária (previously TC3's content)
$ws.Range("B20").Value = "Chefe Clica para realizar o cancelamento de uma diária."
$ws.Range("D20").Value = "SYSTEM Verifica que a solicitação está em situação SOLICITADA; Exibe mensagem de confirmação (MSG987 - Cancelar solicitação de diária) para o usuário (que deve confirmar); Cancela a diária, mudando sua situação para CANCELADA (ver diagrama de estados da diária)."

# TC3's step-2 now describes filtering the listing (previously TC2's content)
$ws.Range("B28").Value = "Chefe Indica alguns parâmetros específicos para a busca; Informa o nome do beneficiário; Filtra a listagem de solicitações."
$ws.Range("D28").Value = "SYSTEM Exibe uma nova listagem de solicitações, de acordo com os filtros informados pelo usuário."
